# Commit: "Optimize the parsing process from two-passe to one and refactor
# the strcuture of table compiler"
#
# Net effect on Assets/ConfigTable/Test.xlsx:
#   - a second worksheet ("Sheet1") is appended after "Test", containing a
#     copy of the same id/name/price/isTest/factor/color table (minus the
#     color legend helper columns I:K) with the id column shifted by +5
#     (the table compiler now also emits/reads a second "page" of rows)
#   - the "Test" sheet's own data is untouched
#   - both sheets end up zoomed to 175% with a new cell selection

$wb  = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("Test")

# --- build the new sheet right after "Test" ------------------------------
$new = $wb.Worksheets.Add($null, $src)
$new.Name = "Sheet1"

# Copy the header + data table (A:F) and the helper rate-factor cells
# (I3:I4) straight from "Test" so styles/number formats come along for
# free; formulas are reattached explicitly afterwards since Copy only
# carries over the cached values.
$src.Range("A1:F12").Copy($new.Range("A1"))
$src.Range("I3:I4").Copy($new.Range("I3"))

$new.Range("E3").Formula   = "=`$I`$4*C3"
$new.Range("E4:E12").Formula = "=`$I`$4*C4"

# Shift every id in the copied table by +5 (1..10 -> 6..15).
for ($r = 3; $r -le 12; $r++) {
    $cell = $new.Cells.Item($r, 1)
    $cell.Value = $cell.Value2 + 5
}

$new.Columns.Item(6).EntireColumn.AutoFit()

# --- view tweaks -----------------------------------------------------------
$new.Activate()
$excel.ActiveWindow.Zoom = 175
$new.Range("D2").Select()

$src.Activate()
$excel.ActiveWindow.Zoom = 175
$src.Range("D10").Select()
